# Initial Time.xlsx - "CA 4.0 files test" commit
#
# Changes applied:
#   1. IT!B2: initial simulation year bumped from 2020 -> 2021.
#   2. About!A7: remove the (accidentally-applied) bold formatting, restoring
#      the cell to the workbook's default (unstyled) look.
#   3. The "About" sheet becomes the active/selected tab (instead of "IT").

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets("About")
$wsIT    = $wb.Worksheets("IT")

# 1. Bump the Initial Time year.
$wsIT.Range("B2").Value = 2021

# 2. Clear the bold font override that had been applied to A7 on the About
#    sheet, bringing it back to the sheet's default formatting.
$wsAbout.Range("A7").Font.Bold = $false

# 3. Make "About" the selected/active sheet (it was "IT" before).
$wsAbout.Activate()
